$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 283.54544
$ws.Range("I2").Value = 196.94118
$ws.Range("K2").Value = 196.94118
$ws.Range("M2").Value = -83.94118

$ws.Range("H88").Value = 924.5
$ws.Range("I88").Value = 362.25
$ws.Range("J88").Value = 1299.3334
$ws.Range("K88").Value = 362.25
$ws.Range("L88").Value = 1299.3334
$ws.Range("M88").Value = 43.75
$ws.Range("N88").Value = -2111.3334

$ws.Range("H91").Value = 924.5
$ws.Range("I91").Value = 362.25
$ws.Range("J91").Value = 1299.3334
$ws.Range("K91").Value = 362.25
$ws.Range("L91").Value = 1299.3334
$ws.Range("M91").Value = 1041.75
$ws.Range("N91").Value = -4107.3334

$ws.Range("H98").Value = 680.5
$ws.Range("I98").Value = 673.86664
$ws.Range("K98").Value = 673.86664
$ws.Range("M98").Value = 824.13336

$ws.Range("H122").Value = 680.5
$ws.Range("I122").Value = 673.86664
$ws.Range("K122").Value = 2021.59992
$ws.Range("M122").Value = 428.4000800000001

$ws.Range("H132").Value = 3384.8845
$ws.Range("I132").Value = 2229.3333
$ws.Range("J132").Value = 8238.200000000001
$ws.Range("K132").Value = 6687.999899999999
$ws.Range("L132").Value = 24714.6
$ws.Range("M132").Value = -4157.999899999999
$ws.Range("N132").Value = -29774.6

$ws.Range("H141").Value = 6445.615
$ws.Range("I141").Value = 6469
$ws.Range("J141").Value = 6367.6665
$ws.Range("K141").Value = 19407
$ws.Range("L141").Value = 19102.9995
$ws.Range("M141").Value = -14227
$ws.Range("N141").Value = -29462.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2385.9092
$ws.Range("I61").Value = 2385.9092
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2385.9092
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2173.9092
$ws.Range("N61").Value = $null

$ws.Range("H74").Value = 1402.4117
$ws.Range("I74").Value = 1402.4117
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1402.4117
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -528.4117000000001
$ws.Range("N74").Value = $null

$ws.Range("H77").Value = 1402.4117
$ws.Range("I77").Value = 1402.4117
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 7012.058500000001
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -2644.058500000001
$ws.Range("N77").Value = $null

$ws.Range("H132").Value = 406.57144
$ws.Range("I132").Value = 307.83334
$ws.Range("J132").Value = 999
$ws.Range("K132").Value = 923.5000200000001
$ws.Range("L132").Value = 2997
$ws.Range("M132").Value = 1606.49998
$ws.Range("N132").Value = -8057

$ws.Range("H136").Value = 2385.9092
$ws.Range("I136").Value = 2385.9092
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 7157.7276
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -4607.7276
$ws.Range("N136").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = $null

$ws.Range("H134").Value = 4041.8667
$ws.Range("I134").Value = 4086.7693
$ws.Range("K134").Value = 12260.3079
$ws.Range("M134").Value = -9725.3079

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2329.7856
$ws.Range("I31").Value = 2080.2222
$ws.Range("J31").Value = 2779
$ws.Range("K31").Value = 2080.2222
$ws.Range("L31").Value = 2779
$ws.Range("M31").Value = -1785.2222
$ws.Range("N31").Value = -3369

$ws.Range("H34").Value = 2329.7856
$ws.Range("I34").Value = 2080.2222
$ws.Range("J34").Value = 2779
$ws.Range("K34").Value = 2080.2222
$ws.Range("L34").Value = 2779
$ws.Range("M34").Value = -1878.2222
$ws.Range("N34").Value = -3183

$ws.Range("H58").Value = 2133.5625
$ws.Range("I58").Value = 2102.6428
$ws.Range("J58").Value = 2350
$ws.Range("K58").Value = 2102.6428
$ws.Range("L58").Value = 2350
$ws.Range("M58").Value = -1899.6428
$ws.Range("N58").Value = -2756

$ws.Range("H86").Value = 7988.6665
$ws.Range("I86").Value = 7476.8335
$ws.Range("J86").Value = 8500.5
$ws.Range("K86").Value = 7476.8335
$ws.Range("L86").Value = 8500.5
$ws.Range("M86").Value = -6353.8335
$ws.Range("N86").Value = -10746.5

$ws.Range("H89").Value = 7988.6665
$ws.Range("I89").Value = 7476.8335
$ws.Range("J89").Value = 8500.5
$ws.Range("K89").Value = 37384.1675
$ws.Range("L89").Value = 42502.5
$ws.Range("M89").Value = -31768.1675
$ws.Range("N89").Value = -53734.5

$ws.Range("H132").Value = 1816.6774
$ws.Range("I132").Value = 1788.75
$ws.Range("K132").Value = 5366.25
$ws.Range("M132").Value = -2836.25

$ws.Range("H134").Value = 1022.1739
$ws.Range("I134").Value = 1022.1739
$ws.Range("K134").Value = 3066.5217
$ws.Range("M134").Value = -531.5217000000002

$ws.Range("H136").Value = 2133.5625
$ws.Range("I136").Value = 2102.6428
$ws.Range("J136").Value = 2350
$ws.Range("K136").Value = 6307.928400000001
$ws.Range("L136").Value = 7050
$ws.Range("M136").Value = -3757.928400000001
$ws.Range("N136").Value = -12150

$ws.Range("H141").Value = 36133.355
$ws.Range("J141").Value = 35452.152
$ws.Range("L141").Value = 35452.152
$ws.Range("N141").Value = -45812.152

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 122883.78
$ws.Range("I2").Value = 84733.16
$ws.Range("J2").Value = 222075.4
$ws.Range("K2").Value = 508398.96
$ws.Range("L2").Value = 1332452.4
$ws.Range("M2").Value = -508285.96
$ws.Range("N2").Value = -1332678.4

$ws.Range("H10").Value = 48.333332
$ws.Range("I10").Value = 22.5
$ws.Range("J10").Value = 100
$ws.Range("K10").Value = 67.5
$ws.Range("L10").Value = 300
$ws.Range("M10").Value = 71.5
$ws.Range("N10").Value = -578

$ws.Range("H47").Value = 1220.6
$ws.Range("I47").Value = 1220.6
$ws.Range("K47").Value = 3661.8
$ws.Range("M47").Value = -3230.8

$ws.Range("H75").Value = 9053.375
$ws.Range("I75").Value = 1856.5
$ws.Range("J75").Value = 11452.333
$ws.Range("K75").Value = 5569.5
$ws.Range("L75").Value = 34356.999
$ws.Range("M75").Value = -4571.5
$ws.Range("N75").Value = -36352.999

$ws.Range("H78").Value = 9053.375
$ws.Range("I78").Value = 1856.5
$ws.Range("J78").Value = 11452.333
$ws.Range("K78").Value = 16708.5
$ws.Range("L78").Value = 103070.997
$ws.Range("M78").Value = -11716.5
$ws.Range("N78").Value = -113054.997

$ws.Range("H87").Value = 4671.6665
$ws.Range("I87").Value = 4000
$ws.Range("J87").Value = 5007.5
$ws.Range("K87").Value = 12000
$ws.Range("L87").Value = 15022.5
$ws.Range("M87").Value = -10752
$ws.Range("N87").Value = -17518.5

$ws.Range("H90").Value = 4671.6665
$ws.Range("I90").Value = 4000
$ws.Range("J90").Value = 5007.5
$ws.Range("K90").Value = 36000
$ws.Range("L90").Value = 45067.5
$ws.Range("M90").Value = -29760
$ws.Range("N90").Value = -57547.5

$ws.Range("H113").Value = 1415.3125
$ws.Range("J113").Value = 1332.1428
$ws.Range("L113").Value = 3996.4284
$ws.Range("N113").Value = -8336.428400000001

$ws.Range("H121").Value = 10257.588
$ws.Range("I121").Value = 19416.166
$ws.Range("J121").Value = 5262
$ws.Range("K121").Value = 58248.49800000001
$ws.Range("L121").Value = 15786
$ws.Range("M121").Value = -56938.49800000001
$ws.Range("N121").Value = -18406

$ws.Range("H140").Value = 11078.75
$ws.Range("I140").Value = 1322
$ws.Range("K140").Value = 3966
$ws.Range("M140").Value = 1214

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 39941.832
$ws.Range("J94").Value = 39941.832
$ws.Range("L94").Value = 39941.832
$ws.Range("N94").Value = -41293.832

$ws.Range("H107").Value = 2054.158
$ws.Range("J107").Value = 2613.6
$ws.Range("L107").Value = 2613.6
$ws.Range("N107").Value = -6453.6

$ws.Range("H113").Value = 1610
$ws.Range("I113").Value = 1332
$ws.Range("K113").Value = 1332
$ws.Range("M113").Value = 838

$ws.Range("H126").Value = 5931.5557
$ws.Range("I126").Value = 6649.5
$ws.Range("K126").Value = 19948.5
$ws.Range("M126").Value = -17478.5

$ws.Range("H132").Value = 1198.1666
$ws.Range("I132").Value = 1261.7273
$ws.Range("J132").Value = 499
$ws.Range("K132").Value = 3785.1819
$ws.Range("L132").Value = 1497
$ws.Range("M132").Value = -1255.1819
$ws.Range("N132").Value = -6557

$ws.Range("H134").Value = 60000
$ws.Range("J134").Value = 60000
$ws.Range("L134").Value = 180000
$ws.Range("N134").Value = -185070

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 617.1
$ws.Range("I22").Value = 607.75
$ws.Range("J22").Value = 631.125
$ws.Range("K22").Value = 607.75
$ws.Range("L22").Value = 631.125
$ws.Range("M22").Value = -312.75
$ws.Range("N22").Value = -1221.125

$ws.Range("H27").Value = 617.1
$ws.Range("I27").Value = 607.75
$ws.Range("J27").Value = 631.125
$ws.Range("K27").Value = 607.75
$ws.Range("L27").Value = 631.125
$ws.Range("M27").Value = -500.75
$ws.Range("N27").Value = -845.125

$ws.Range("H61").Value = 4016.3333
$ws.Range("I61").Value = 3774.5
$ws.Range("K61").Value = 3774.5
$ws.Range("M61").Value = -3572.5

$ws.Range("H113").Value = 4016.3333
$ws.Range("I113").Value = 3774.5
$ws.Range("K113").Value = 3774.5
$ws.Range("M113").Value = -1604.5

$ws.Range("H122").Value = 5911.5
$ws.Range("I122").Value = 5563.5415
$ws.Range("K122").Value = 16690.6245
$ws.Range("M122").Value = -14240.6245

$ws.Range("H132").Value = 1774.1578
$ws.Range("I132").Value = 1260.6666
$ws.Range("K132").Value = 3781.9998
$ws.Range("M132").Value = -1251.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1138.0769
$ws.Range("I107").Value = 1257.4286
$ws.Range("K107").Value = 3772.2858
$ws.Range("M107").Value = -1852.2858

$ws.Range("H122").Value = 1671
$ws.Range("I122").Value = 1326.7142
$ws.Range("J122").Value = 2273.5
$ws.Range("K122").Value = 3980.1426
$ws.Range("L122").Value = 6820.5
$ws.Range("M122").Value = -1530.1426
$ws.Range("N122").Value = -11720.5

$ws.Range("H136").Value = 2894.2693
$ws.Range("I136").Value = 2894.2693
$ws.Range("K136").Value = 8682.8079
$ws.Range("M136").Value = -6132.8079
